$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update RE 9 (row 11) description: now mentions "menores preços" and "em um período"
$ws.Range("B11").Value2 = "Disponibilizar histórico e variação dos menores preços encontrados em um período para os produtos da cesta (própria ou do sistema)"

# 2. Make room for the new "RE 14" use case row by shifting the trailing
#    "A pensar" / "Imaginar..." notes rows down by one.
$ws.Rows("17:17").Insert()

# 3. Fill in the new requirement row (row 16), which previously only had
#    the RE number (A16=14) populated.
$ws.Range("B16").Value2 = "Consultar o local mais barato onde se pode adquirir um produto ou todos os produtos da cesta na data atual."
$ws.Range("C16").Value2 = "Consumidor"
$ws.Range("E16").Value2 = 13

$ws.Range("A16").VerticalAlignment = -4160
$ws.Range("B16").WrapText = $true
$ws.Rows("16:16").RowHeight = 30

# 4. Append a new "A pensar" bullet about allowing basket quantities.
$ws.Range("B20").Value2 = "Permitir montar a cesta incluindo quantidades dos produtos (na versão inicial, apenas uma unidade de cada produto)"
$ws.Range("B20").WrapText = $true
$ws.Rows("20:20").RowHeight = 30

# F16 is filled in last (matches the order new shared strings were appended).
$ws.Range("F16").Value2 = "Consultar Menor Preço de Cesta"

# 5. Selection recorded in the saved workbook (also clears the old
#    "topLeftCell" scroll position left over from the previous edit).
$ws.Range("F14").Select()
